$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: Row index, Column index, old text (for a sanity check),
# new text. Row 1 is the header row; data rows are 2-13.
#
# NOTE: Find.Execute on this runtime performs a document-wide
# replace (it ignores the Range it is invoked on for scoping purposes),
# so using it per-cell would bleed across cells whenever one value's text
# is a substring of another (e.g. "5.5" also occurs inside "7 (5.6)"/
# "12 (25.5)" once other edits are applied). Instead we address each
# table cell directly via Rows/Cells and assign Range.Text, which only
# touches that cell and keeps its run formatting intact.
$edits = @(
    @(2, 5, "11 ( 8.7)", "12 ( 9.4)"),
    @(2, 7, "7 (5.6)",   "7 (5.5)"),
    @(3, 5, "27 (21.4)", "28 (22.0)"),
    @(4, 5, "40 (31.7)", "41 (32.3)"),
    @(5, 3, "5.5",       "5.0"),
    @(5, 5, "51 (40.5)", "52 (40.9)"),
    @(5, 7, "6 (4.8)",   "6 (4.7)"),
    @(6, 5, "55 (43.7)", "56 (44.1)"),
    @(7, 3, "2.5",       "2.0"),
    @(7, 5, "61 (48.4)", "62 (48.8)"),
    @(8, 5, "68 (54.0)", "69 (54.3)"),
    @(9, 5, "66 (52.4)", "67 (52.8)"),
    @(10, 5, "71 (56.3)", "72 (56.7)"),
    @(11, 5, "74 (58.7)", "75 (59.1)"),
    @(12, 5, "76 (60.3)", "77 (60.6)"),
    @(13, 5, "81 (64.3)", "82 (64.6)")
)

foreach ($edit in $edits) {
    $rowIdx = $edit[0]
    $colIdx = $edit[1]
    $oldText = $edit[2]
    $newText = $edit[3]

    $cell = $t.Rows.Item($rowIdx).Cells.Item($colIdx)
    $range = $cell.Range
    # Cell.Range.Text includes the trailing cell-mark characters
    # (paragraph mark + cell mark), so trim those before comparing.
    $current = $range.Text.TrimEnd([char]13, [char]7)
    if ($current -eq $oldText) {
        $range.Text = $newText
    }
}
